# Insert a new data row at row 812 (pushes existing rows 812:924 down to 813:925)
# and populate it with a new "Camote" price record ("1a (guarda)").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(812).Insert()

$ws.Cells.Item(812, 1).Value = 8
$ws.Cells.Item(812, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(812, 3).Value = "Coquimbo"
$ws.Cells.Item(812, 4).Value = 45077
$ws.Cells.Item(812, 5).Value = 4
$ws.Cells.Item(812, 6).Value = 100112045
$ws.Cells.Item(812, 7).Value = "Zapallo"
$ws.Cells.Item(812, 8).Value = "Camote"
$ws.Cells.Item(812, 9).Value = "1a (guarda)"
$ws.Cells.Item(812, 10).Value = 1400
$ws.Cells.Item(812, 11).Value = 400
$ws.Cells.Item(812, 12).Value = 500
$ws.Cells.Item(812, 13).Value = 450
$ws.Cells.Item(812, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(812, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(812, 16).Value = 450
$ws.Cells.Item(812, 17).Value = 1
$ws.Cells.Item(812, 18).Value = "Hortaliza"
